# Normalize "Recorded By" (column G) values: the list of recorder names/emails
# in each cell was stored in reverse order and needs to be flipped back,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Rows whose value contains "admin@admin.com" are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G" + $r)
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }
    if (-not ($val -is [string])) {
        continue
    }
    if ($val.IndexOf(",") -lt 0) {
        continue
    }
    if ($val.Contains("admin@admin.com")) {
        continue
    }

    $parts = $val.Split(",")
    $n = $parts.Length
    $trimmed = @()
    for ($i = 0; $i -lt $n; $i++) {
        $trimmed += $parts[$i].Trim()
    }
    $revIdx = ($n - 1)..0
    $rev = $trimmed[$revIdx]
    $joined = [string]::Join(", ", $rev)

    $cell.Value = $joined
}
